$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: "Objetivos:" label stays; fix body text (B10/C10) ---
$ws.Range("B10").Value = 'A disciplina visa apresentar aos estudantes o conceito, tipos, modelos e sistemas de eco-inovação para o desenvolvimento da capacidade analítica e propositiva como competências profissionais nas áreas de inovação e sustentabilidade.'
$ws.Range("C10").Value = 'A disciplina visa apresentar aos estudantes o conceito, tipos, modelos e sistemas de eco-inovação para o desenvolvimento da capacidade analítica e propositiva como competências profissionais nas áreas de inovação e sustentabilidade.'

# --- Row 12: was "Docentes responsaveis:" (label only, no B/C) ---
# --- becomes "Programa resumido:" with Gustavo-name body, height 60 ---
# Copy B11:C11 formatting first so the brand-new B12/C12 cells inherit
# the correct column styles (wrap-text / red-font) instead of defaulting.
$ws.Range("B11:C11").Copy()
$ws.Range("B12:C12").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A12").Value = 'Programa resumido:'
$ws.Range("B12").Value = '5840820 - Gustavo Aristides Santana Martinez'
$ws.Range("C12").Value = '5840820 - Gustavo Aristides Santana Martinez'
$ws.Rows.Item(12).RowHeight = 60

# --- Row 13: was "Programa resumido:"/"Semestral" -> "Short syllabus:" + EN short syllabus ---
$ws.Range("A13").Value = 'Short syllabus:'
$ws.Range("B13").Value = 'Eco-innovation. Eco-innovation metrics. Introduction to products life-cycle. Eco-innovation in the industry.  Case study of Eco-innovation projects in Brazil. Methods and tools to support the process of Eco-innovation. Early identification of failure as support to Eco-innovation. TRIZ as a response to Eco-innovation. Methodological proposal for Eco-innovative solutions.'
$ws.Range("C13").Value = 'Eco-innovation. Eco-innovation metrics. Introduction to products life-cycle. Eco-innovation in the industry.  Case study of Eco-innovation projects in Brazil. Methods and tools to support the process of Eco-innovation. Early identification of failure as support to Eco-innovation. TRIZ as a response to Eco-innovation. Methodological proposal for Eco-innovative solutions.'

# --- Row 14: was "Short syllabus:"/EN short syllabus -> "Programa:" + NEW PT full programa text, ht 60->120 ---
$ws.Range("A14").Value = 'Programa:'
$ws.Range("B14").Value = 'Eco inovação. Métricas da eco-inovação. Introdução ao Ciclo de vida do produto. Eco inovação na indústria. Estudo de casos de projetos de eco-inovação no Brasil. Métodos e ferramentas suporte do processo de eco-inovação. Identificação antecipada de falha como suporte a eco-inovação. TRIZ como resposta a eco-inovação. Proposta metodológica para soluções eco-inovadoras.'
$ws.Range("C14").Value = 'Eco inovação. Métricas da eco-inovação. Introdução ao Ciclo de vida do produto. Eco inovação na indústria. Estudo de casos de projetos de eco-inovação no Brasil. Métodos e ferramentas suporte do processo de eco-inovação. Identificação antecipada de falha como suporte a eco-inovação. TRIZ como resposta a eco-inovação. Proposta metodológica para soluções eco-inovadoras.'
$ws.Rows.Item(14).RowHeight = 120

# --- Row 15: was "Programa:"/"01/01/2020" -> "Syllabus:" + EN full syllabus text, stays ht 120 ---
$ws.Range("A15").Value = 'Syllabus:'
$ws.Range("B15").Value = '1. Eco-innovation: concepts, determinant factors, barriers, types of Eco-innovative agents, category of Eco-innovations. 2.  Eco-innovation metrics: Andersen metrics, Arundel & Kemp metrics, OECD metrics.3. Introduction to products life-cycle: analytical perspective, production chain analysis, Green Supply Chain Management Practices.4.  Eco-innovation in the industry: chemistry, agro-food, metal mechanics.5.  Case study of Eco-innovation projects in Brazil.6. Methods and tools to support the process of Eco-innovation: Eco-Compass, Eco-Ideation Tool, Value Mapping Tool, Design for Environment (DfE), EcoASIT, others.7. Early identification of failure as support to Eco-innovation: problem, scene, resources.8. TRIZ as a response to Eco-innovation: inventive principles, engineering parameters, contradictions matrix.9. Methodological proposal for Eco-innovative solutions in technological categories: definition, measurement, analysis, creation'
$ws.Range("C15").Value = '1. Eco-innovation: concepts, determinant factors, barriers, types of Eco-innovative agents, category of Eco-innovations. 2.  Eco-innovation metrics: Andersen metrics, Arundel & Kemp metrics, OECD metrics.3. Introduction to products life-cycle: analytical perspective, production chain analysis, Green Supply Chain Management Practices.4.  Eco-innovation in the industry: chemistry, agro-food, metal mechanics.5.  Case study of Eco-innovation projects in Brazil.6. Methods and tools to support the process of Eco-innovation: Eco-Compass, Eco-Ideation Tool, Value Mapping Tool, Design for Environment (DfE), EcoASIT, others.7. Early identification of failure as support to Eco-innovation: problem, scene, resources.8. TRIZ as a response to Eco-innovation: inventive principles, engineering parameters, contradictions matrix.9. Methodological proposal for Eco-innovative solutions in technological categories: definition, measurement, analysis, creation'

# --- Row 16: was "Syllabus:"/EN full syllabus -> "Avaliacao:" (label only), ht 120 -> none ---
$ws.Range("A16").Value = 'Avaliação:'
$ws.Range("B16:C16").Clear()
$ws.Rows.Item(16).AutoFit()

# --- Row 17: was "Avaliacao:" (label only) -> "Metodo:" + NEW PT full syllabus text, ht none->60 ---
$ws.Range("B10:C10").Copy()
$ws.Range("B17:C17").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A17").Value = 'Método:'
$ws.Range("B17").Value = '1. Eco-inovação: conceitos, fatores determinantes, barreiras, tipos de agentes eco-inovadores, categorias de eco inovações.2. Métricas da eco-inovação: métricas de Andersen, métricas de Arundel & Kemp, métricas da OECD.3. Introdução ao Ciclo de vida do produto: perspectiva analítica, análise da cadeia de produção, práticas de Green Supply Chain Management.4. Eco inovação na indústria: química, agro alimentos, metal mecânica.5. Estudo de casos de projetos de eco inovação no Brasil.6. Métodos e ferramentas suporte do processo de eco-inovação: Eco-Compass, Eco-Ideation Tool, Value Mapping Tool, Design for Environment (DfE), EcoASIT, outros.7. Identificação antecipada de falha como suporte a eco-inovação: o problema, o cenário, os recursos.8. TRIZ como resposta a eco inovação: princípios inventivos, parâmetros de engenharia, matriz das contradições.9. Proposta metodológica para soluções eco inovadoras na categoria tecnologias: definir, medir, analisar, criar'
$ws.Range("C17").Value = '1. Eco-inovação: conceitos, fatores determinantes, barreiras, tipos de agentes eco-inovadores, categorias de eco inovações.2. Métricas da eco-inovação: métricas de Andersen, métricas de Arundel & Kemp, métricas da OECD.3. Introdução ao Ciclo de vida do produto: perspectiva analítica, análise da cadeia de produção, práticas de Green Supply Chain Management.4. Eco inovação na indústria: química, agro alimentos, metal mecânica.5. Estudo de casos de projetos de eco inovação no Brasil.6. Métodos e ferramentas suporte do processo de eco-inovação: Eco-Compass, Eco-Ideation Tool, Value Mapping Tool, Design for Environment (DfE), EcoASIT, outros.7. Identificação antecipada de falha como suporte a eco-inovação: o problema, o cenário, os recursos.8. TRIZ como resposta a eco inovação: princípios inventivos, parâmetros de engenharia, matriz das contradições.9. Proposta metodológica para soluções eco inovadoras na categoria tecnologias: definir, medir, analisar, criar'
$ws.Rows.Item(17).RowHeight = 60

# --- Row 18: was "Metodo:"/Gustavo name -> "Criterio:" + Aulas expositivas text, stays ht 60 ---
$ws.Range("A18").Value = 'Critério:'
$ws.Range("B18").Value = 'Aulas expositivas, discussão de casos em sala de aula, painéis, debates, seminários, análise de vídeos e palestrantes externos.'
$ws.Range("C18").Value = 'Aulas expositivas, discussão de casos em sala de aula, painéis, debates, seminários, análise de vídeos e palestrantes externos.'

# --- Row 19: was "Criterio:"/Aulas expositivas -> "Norma de recuperacao:" + NF formula text, stays ht 60 ---
$ws.Range("A19").Value = 'Norma de recuperação:'
$ws.Range("B19").Value = 'NF= (N1 + N2)/2Onde: NF = nota final; N = nota'
$ws.Range("C19").Value = 'NF= (N1 + N2)/2Onde: NF = nota final; N = nota'

# --- Row 20: was "Norma de recuperacao:"/NF formula text (ht 60) -> "Bibliografia:" + Estara apto text, ht 60->120 ---
$ws.Range("A20").Value = 'Bibliografia:'
$ws.Range("B20").Value = 'Estará apto a efetuar a prova de reavaliação o aluno que tiver como média final na disciplina uma nota igual ou superior a três (3,0) e inferior a cinco (5,0), e tiver, no mínimo, 70% (setenta por cento) de frequência às aulas. O cálculo de uma média aritmética simples será feito com a nota da prova de reavaliação e a média final obtida pelo aluno na disciplina. Se esta média resultar em nota igual ou superior a cinco (5,0), o aluno será aprovado.'
$ws.Range("C20").Value = 'Estará apto a efetuar a prova de reavaliação o aluno que tiver como média final na disciplina uma nota igual ou superior a três (3,0) e inferior a cinco (5,0), e tiver, no mínimo, 70% (setenta por cento) de frequência às aulas. O cálculo de uma média aritmética simples será feito com a nota da prova de reavaliação e a média final obtida pelo aluno na disciplina. Se esta média resultar em nota igual ou superior a cinco (5,0), o aluno será aprovado.'
$ws.Rows.Item(20).RowHeight = 120

# --- Row 21: was "Bibliografia:"/Estara apto text (ht 120) -> removed (content now lives in row 20) ---
$ws.Rows.Item(21).Delete()

